# Integrate the new minigame ("连连看" / Lianliankan match game) as row 10
# of the "mini" sheet's table, update the theme's light-1 color, and move
# the active selection to the newly added cell — mirroring the authored
# commit "basically integrade the linegame".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data lives in an Excel Table (ListObject) spanning A3:E9. Adding a
# ListRow grows the table (and its AutoFilter range) to A3:E10 and keeps
# the worksheet `dimension` in sync automatically.
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# Populate the new record (row 10): Id, Name, WindowId, IconPath.
# (BgImage / column E is intentionally left blank, matching the source data.)
$ws.Cells.Item(10, 1).Value = 17000007
$ws.Cells.Item(10, 2).Value = "连连看"
$ws.Cells.Item(10, 3).Value = 1106
$ws.Cells.Item(10, 4).Value = "GameButton6"

# Lighten the theme's "Background 1" (lt1) color from the pale green
# CAEACD to pure white FFFFFF.
$scheme = $wb.Theme.ThemeColorScheme
$bg1 = $scheme.Colors(2)
$bg1.RGB = 16777215

# Move/record the active selection on the newly added cell.
$ws.Range("E10").Select()
